$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the changed cells so the numeric-looking values
# are stored as text (matching the source workbook, which holds these
# price/volume columns as plain strings rather than numbers).
$changedRefs = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "E18", "D19", "E19", "D20", "E20", "D21", "E21", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "D26", "E26", "D27", "E27", "D39", "E39", "D40", "E40", "D41", "E41", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "E47", "E48", "D49", "E49", "D50", "E50", "D51", "E51")
foreach ($ref in $changedRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

# Updated coin price / 1h-volume values (symbol list refresh).
$ws.Range("D2").Value = "320.34"
$ws.Range("E2").Value = "4.87%"
$ws.Range("D3").Value = "36.07"
$ws.Range("E3").Value = "-0.23%"
$ws.Range("D4").Value = "5.132"
$ws.Range("E4").Value = "1.17%"
$ws.Range("D5").Value = "0.08163"
$ws.Range("E5").Value = "3.85%"
$ws.Range("D6").Value = "2.145"
$ws.Range("E6").Value = "-2.40%"
$ws.Range("D7").Value = "8.040"
$ws.Range("E7").Value = "1.57%"
$ws.Range("D8").Value = "4.137"
$ws.Range("E8").Value = "1.04%"
$ws.Range("D9").Value = "0.9261"
$ws.Range("E9").Value = "0.60%"
$ws.Range("D10").Value = "0.1007"
$ws.Range("E10").Value = "4.54%"
$ws.Range("D11").Value = "0.1889"
$ws.Range("E11").Value = "1.55%"
$ws.Range("D12").Value = "0.09263"
$ws.Range("E12").Value = "7.77%"
$ws.Range("E13").Value = "2.94%"
$ws.Range("D14").Value = "0.09927"
$ws.Range("E14").Value = "0.04%"
$ws.Range("D15").Value = "0.001435"
$ws.Range("E15").Value = "0.42%"
$ws.Range("D16").Value = "0.005680"
$ws.Range("E16").Value = "0.56%"
$ws.Range("D17").Value = "3.448"
$ws.Range("E17").Value = "-0.46%"
$ws.Range("E18").Value = "15.88%"
$ws.Range("D19").Value = "0.3372"
$ws.Range("E19").Value = "-1.56%"
$ws.Range("D20").Value = "0.1329"
$ws.Range("E20").Value = "-1.94%"
$ws.Range("D21").Value = "5.047"
$ws.Range("E21").Value = "4.92%"
$ws.Range("E22").Value = "-0.87%"
$ws.Range("D23").Value = "0.04584"
$ws.Range("E23").Value = "0.99%"
$ws.Range("D24").Value = "0.001242"
$ws.Range("E24").Value = "0.78%"
$ws.Range("D25").Value = "0.004730"
$ws.Range("E25").Value = "-7.20%"
$ws.Range("D26").Value = "0.0001299"
$ws.Range("E26").Value = "-7.34%"
$ws.Range("D27").Value = "0.0004498"
$ws.Range("E27").Value = "-5.32%"
$ws.Range("D39").Value = "0.02025"
$ws.Range("E39").Value = "11.05%"
$ws.Range("D40").Value = "0.05001"
$ws.Range("E40").Value = "4.58%"
$ws.Range("D41").Value = "0.007814"
$ws.Range("E41").Value = "0.22%"
$ws.Range("E42").Value = "0.20%"
$ws.Range("D43").Value = "0.007654"
$ws.Range("E43").Value = "-1.10%"
$ws.Range("D44").Value = "0.002094"
$ws.Range("E44").Value = "-6.21%"
$ws.Range("D45").Value = "0.01219"
$ws.Range("E45").Value = "8.58%"
$ws.Range("D46").Value = "0.00006478"
$ws.Range("E46").Value = "4.45%"
$ws.Range("E47").Value = "-0.14%"
$ws.Range("E48").Value = "17.43%"
$ws.Range("D49").Value = "0.001899"
$ws.Range("E49").Value = "-5.05%"
$ws.Range("D50").Value = "0.00002099"
$ws.Range("E50").Value = "-0.14%"
$ws.Range("D51").Value = "0.0001999"
$ws.Range("E51").Value = "-0.14%"
